$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D: sequential numbering 1..36 for rows 6-41 (new "Sl. No." style column
# added alongside the existing Topic/Problem/Done columns).
for ($row = 6; $row -le 41; $row++) {
    $ws.Cells.Item($row, 4).Value = $row - 5
}

# Row 80: this question's status was unclear -> flag it in red and
# change the "Done" column text to "not clear".
$ws.Range("C80").Value = "not clear"
$ws.Range("B80").Font.Color = 255

# Rows 85-98: these questions are now finished -> flag them in green and
# change the "Done" column text to "yes".
for ($row = 85; $row -le 98; $row++) {
    $ws.Range("C$row").Value = "yes"
    $ws.Range("B$row").Font.Color = 0x50B000
}

# Scroll the sheet down and move the active selection to B98, reflecting
# where work left off.
[void]$ws.Range("B98").Select()
$excel.ActiveWindow.ScrollRow = 82
$excel.ActiveWindow.ScrollColumn = 1
